# "Drop in all data files from 3.0 RMI script"
#
# This workbook edit:
#  1. Deletes the "Texas Notes" worksheet entirely (it was specific to an
#     earlier, Texas-only version of the model and is no longer needed).
#  2. Updates the DR discount rate value on the "DR" sheet (cell B2) from
#     5.87% down to a flat 3%.
#  3. Leaves the saved selection on each remaining sheet matching the
#     final on-screen state (About!A16:A18 selected, with "About" as the
#     active tab; DR!B1 selected on the DR sheet).

$wb = $excel.ActiveWorkbook

# Avoid any "are you sure you want to delete this sheet" prompt.
$excel.DisplayAlerts = $false

# 1. Remove the "Texas Notes" worksheet.
$texasNotes = $wb.Worksheets.Item("Texas Notes")
$texasNotes.Delete()

# 2. Update the discount rate value on the DR sheet: 5.87% -> 3%.
$dr = $wb.Worksheets.Item("DR")
$dr.Range("B2").Value = 0.03

# 3. Restore the selections/active sheet recorded in the workbook view.
$about = $wb.Worksheets.Item("About")
$about.Activate()
$about.Range("A16:A18").Select()

$dr.Activate()
$dr.Range("B1").Select()

# "About" remains the selected/active tab in the saved workbook.
$about.Activate()

$excel.DisplayAlerts = $true
